$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# "Elimina EC anteriores y se agregan nuevos, se modifica base de datos":
# the account-statement detail rows (16-23) are refreshed with the new
# "base de datos" ordering. The "Periodo Mora" (col E) / "Valor Mora"
# (col F) pairs for that block end up reversed top to bottom:
#   before: 2010, 2011, 2012, 2101, 2102, 2103, 2104, 2105
#           (Valor Mora 35112 for all but the last period, which is 28090)
#   after:  2105, 2104, 2103, 2102, 2101, 2012, 2011, 2010
#           (Valor Mora 28090 for the first period now, 35112 for the rest)
$periods = @("2105", "2104", "2103", "2102", "2101", "2012", "2011", "2010")
$amounts = @(28090, 35112, 35112, 35112, 35112, 35112, 35112, 35112)

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]
    $ws.Cells.Item($row, 6).Value = $amounts[$i]
}
